$d = $word.ActiveDocument

# --- Step 1: force Word to mint a fresh numbered-list definition (decimal,
# 9 levels) in numbering.xml -- this is what a real user gets from the
# Numbering-Library "1. 2. 3." gallery entry. We do this against a throwaway
# paragraph appended at the very end of the story, then delete that
# paragraph again so the only lasting effect is the new abstractNum/num
# pair Word minted (which InsertXML below can then reference via numId=2).
$tailCount = $d.Paragraphs.Count
$tailPara = $d.Paragraphs($tailCount)
$tailPara.Range.InsertParagraphAfter()
$scratchCount = $d.Paragraphs.Count
$scratch = $d.Paragraphs($scratchCount)
$scratch.Range.Text = "scratch"
$scratch.Range.ListFormat.ApplyNumberDefault()
$scratch.Range.Delete()

# --- Step 2: locate the "3.3 Scope" paragraph that ends in "...Python-based
# tools." (the second-to-last paragraph in the original document -- the
# last paragraph is the trailing blank one before the sectPr) and insert a
# fresh empty paragraph right after it. That new, empty paragraph is then
# replaced in one shot via InsertXML with the full "3.4 Project Plan"
# section (heading, intro sentence, six numbered list items, and the
# closing paragraph) so every run/proofErr/numPr comes out exactly as
# authored.
$n = $d.Paragraphs.Count
$scopePara = $d.Paragraphs($n - 1)
$scopePara.Range.InsertParagraphAfter()
$placeholder = $d.Paragraphs($n)
$target = $placeholder.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body><w:p/><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>3.4 Project Plan</w:t></w:r></w:p><w:p><w:r><w:t>The project follows a structured machine learning workflow inspired by CRISP-DM:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Data Collection &amp; Understanding</w:t></w:r><w:r><w:t>: Load and explore the EU housing dataset, perform statistical summaries, and visualize trends.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Data Preparation</w:t></w:r><w:r><w:t>: Clean data, encode categorical features, scale numeric values, and engineer new features (e.g., ratios and interactions).</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Modeling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: Train and tune six models — Linear Regression, Random Forest, Gradient Boosting, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>XGBoost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, SVR, and KNN — using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GridSearchCV</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for optimization.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Evaluation</w:t></w:r><w:r><w:t>: Use R², MAE, and RMSE to assess performance. Visualize predictions and residuals.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Comparison &amp; Insights</w:t></w:r><w:r><w:t>: Determine the best-performing model and identify key predictive features driving housing prices.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Conclusion</w:t></w:r><w:r><w:t xml:space="preserve">: Summarize findings and propose future work, including time-series </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>modeling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> or integration of external data sources.</w:t></w:r></w:p><w:p><w:r><w:t>By combining time-aware feature engineering with non-linear models and explainable AI techniques, this project aims to deliver a robust, interpretable solution for housing price forecasting in the EU (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Gyourko</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> et al., 2021; Kaufmann and Steinmetz, 2020).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
